$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("summary")

# Fill in row 22 with the new diary entry: date, begin, end, (duration formula
# already present), remarks
$ws.Range("A22").NumberFormat = "@"
$ws.Range("A22").Value = "12.11.18"
$ws.Range("A22").NumberFormat = "[$-F800]dddd\,\ mmmm\ dd\,\ yyyy"
$ws.Range("B22").Value = 0.4375
$ws.Range("C22").Value = 0.75
$ws.Range("E22").Value = "-Output Options"

# Move the active selection to A23 (next empty row), matching the saved
# workbook view state after the edit
$ws.Range("A23").Select()
